# Pushing Extent Report changes
#
# The Configuration sheet tracks the current Chrome version under
# "Current Version" (column E) next to "Chrome" (column D). Bump it to
# the latest build.

$wb = $excel.ActiveWorkbook
$configSheet = $wb.Worksheets.Item("Configuration")

$configSheet.Range("E2").Value = "103.0.5060.114 "

# Activate the Configuration sheet and move the selection, matching the
# author's last recorded cursor position when the change was saved.
$configSheet.Activate()
$configSheet.Range("C11").Select()
